$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (the old header row) is no longer used - clear it out entirely.
$ws.Rows.Item(1).ClearContents()

# Rows 2-4 had "test" in column A; it should read "Name" like the rest.
$ws.Range("A2").Value = "Name"
$ws.Range("A3").Value = "Name"
$ws.Range("A4").Value = "Name"

# Column C (the old "Price" column) gets new values.
$ws.Range("C2").Value = "8,9"
$ws.Range("C3").Value = "8,9"

# "8.8" and "10.0" look like numbers, so force the cells to text first
# (otherwise Excel auto-converts them to numeric values), then drop the
# temporary number format back to the default so no stray formatting
# sticks around on the cell.
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "8.8"
$ws.Range("C4").Style = "Normal"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "10.0"
$ws.Range("C5").Style = "Normal"

# Rows 6-10 and column D are no longer part of the data - remove them.
$ws.Range("A6:D10").ClearContents()
